$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.246.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.859.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7113'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9992'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07968'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.97%  '
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08182'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.829.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.178'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7046'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.99%  '
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.212.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.859'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007879'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9981'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.082.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9993'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.424'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.942'
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.934'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.434'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.477'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.377'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.022'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05227'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7116'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9991'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.663'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.727'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9293'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.130.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4276'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.842'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9989'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("E48").Value = '  -4.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.770'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.165'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.973.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.81%  '
